$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header "molar_flow_in" in I1
$ws.Range("I1").ClearContents()

# Clear the value in I2 but keep its number formatting
$ws.Range("I2").ClearContents()

# Update the selection to H3 (matches the selection change in the diff)
$ws.Range("H3").Select()
